$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.585.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.154.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.60"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.30%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.80"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.472.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.46%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.44%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.162.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.530.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.16%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.70%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.58"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.85%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.66%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +8.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +16.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.95%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.60"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.511.90"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0920"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.99%  "

$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000189"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +26.81%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.09%  "

